$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@('000315', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_003_SavedCycle_1_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_003_SavedCycle_1_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:24')
  ,@('000316', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_003_SavedCycle_2_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_003_SavedCycle_2_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:25')
  ,@('000317', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_003_SavedCycle_3_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_003_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:25')
  ,@('000318', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_003_SavedCycle_4_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_003_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:25')
  ,@('000319', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_003_SavedCycle_5_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_003_SavedCycle_5_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:25')
  ,@('000320', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_004_SavedCycle_1_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_004_SavedCycle_1_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:25')
  ,@('000321', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_004_SavedCycle_2_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_004_SavedCycle_2_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:25')
  ,@('000322', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_004_SavedCycle_3_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_004_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:25')
  ,@('000323', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_004_SavedCycle_4_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_004_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:25')
  ,@('000324', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_004_SavedCycle_5_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_004_SavedCycle_5_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:25')
  ,@('000325', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_005_SavedCycle_1_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_005_SavedCycle_1_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:26')
  ,@('000326', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_005_SavedCycle_2_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_005_SavedCycle_2_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:26')
  ,@('000327', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_005_SavedCycle_3_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_005_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:26')
  ,@('000328', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_005_SavedCycle_4_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_005_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:26')
  ,@('000329', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_005_SavedCycle_5_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_005_SavedCycle_5_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:26')
  ,@('000330', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_006_SavedCycle_1_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_006_SavedCycle_1_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:26')
  ,@('000331', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_006_SavedCycle_2_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_006_SavedCycle_2_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:26')
  ,@('000332', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_006_SavedCycle_3_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_006_SavedCycle_3_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:26')
  ,@('000333', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_006_SavedCycle_4_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_006_SavedCycle_4_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:26')
  ,@('000334', 'Details: 60bpm_mf_44path_Normal_SplitCycles_TimeBased_Session01_ExpA1_006_SavedCycle_5_Resampled. Script used: ResampledVisualisationCycleAverage_XY.  Dataset used: C:\Users\Courtney\source\repos\ThesisProject\Data\Session01_ManipulatedData\SavedCycles_Resampled\Session01_ExpA1_006_SavedCycle_5_Resampled.mat. File Location: Visualisations/Session01_Resampled_Average_Figures. Date Generated: 12-May-2023 15:27:26')
  ,@('000335', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath.. Script used: PcRegisterIcpTest.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_006_SavedCycle_4_Resampled. File Location: Visualisations/RegistrationAnalysis. Date Generated: 12-May-2023 22:15:23')
  ,@('000336', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath.. Script used: PcRegisterIcpTest.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_006_SavedCycle_4_Resampled. File Location: Visualisations/RegistrationAnalysis. Date Generated: 12-May-2023 22:16:24')
  ,@('000337', 'Details: Rigid Registration + Scaling - Session01_ManipulatedData - Euclidean distance underneath.. Script used: PcRegisterIcpTest.  Dataset used: reference: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_All_Resampled_Average. data: Session01_ManipulatedData/SavedCycles_Resampled/Session01_ExpA1_006_SavedCycle_4_Resampled. File Location: Visualisations/RegistrationAnalysis. Date Generated: 12-May-2023 22:16:28')
)

$startRow = 316
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $guid = $data[$i][0]
    $details = $data[$i][1]

    $aCell = $ws.Cells.Item($row, 1)
    $aCell.NumberFormat = "@"
    $aCell.Value = $guid
    $aCell.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $details
}

Write-Output ("Updated range: " + $ws.UsedRange.Address())
